# Apply updated crypto price/volume figures (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.148.24"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.59"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("E5").Value = "  -4.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5091"
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2580"
$ws.Range("E8").Value = "  -3.32%  "
$ws.Range("E9").Value = "  -4.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.91"
$ws.Range("E10").Value = "  -4.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07806"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.649.13"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.278"
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.882.21"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5504"
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8005"
$ws.Range("E16").Value = "  -2.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.92"
$ws.Range("E17").Value = "  -6.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.159.87"
$ws.Range("E18").Value = "  -4.40%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.68"
$ws.Range("E20").Value = "  -6.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.404"
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.04"
$ws.Range("E22").Value = "  -3.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.002"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.75"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.736"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("E27").Value = "  -3.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.965"
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.78"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05121"
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.241"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.347"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.216"
$ws.Range("E33").Value = "  -6.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.566"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.749"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.373"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9284"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5692"
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.155.65"
$ws.Range("E39").Value = "  +5.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01590"
$ws.Range("E40").Value = "  -3.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.560"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8340"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.638"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.28"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.791.95"
$ws.Range("E46").Value = "  -3.37%  "
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.63"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.889"
$ws.Range("E51").Value = "  -2.97%  "
